$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 13.92120970214492
$ws.Cells.Item(2, 2).Value = 10.96131004207145
$ws.Cells.Item(2, 3).Value = 18.26922699981088
$ws.Cells.Item(3, 1).Value = 18.41327381060561
$ws.Cells.Item(3, 2).Value = 6.939884531963467
$ws.Cells.Item(3, 3).Value = 29.82095882800603
$ws.Cells.Item(4, 1).Value = 17.62925599830132
$ws.Cells.Item(4, 2).Value = 12.57160958715028
$ws.Cells.Item(4, 3).Value = 23.26916816943996
$ws.Cells.Item(5, 1).Value = 22.07706730589455
$ws.Cells.Item(5, 2).Value = 8.598091636768849
$ws.Cells.Item(5, 3).Value = 34.93086206240482
$ws.Cells.Item(6, 1).Value = 22.47389071713425
$ws.Cells.Item(6, 2).Value = 18.87719289283462
$ws.Cells.Item(6, 3).Value = 26.14242433608779
$ws.Cells.Item(7, 1).Value = 12.70541138739235
$ws.Cells.Item(7, 2).Value = 10.23569401058053
$ws.Cells.Item(7, 3).Value = 16.0334719875313
$ws.Cells.Item(8, 1).Value = 1.340504840444483
$ws.Cells.Item(8, 2).Value = 1.060471962944
$ws.Cells.Item(8, 3).Value = 1.845320729711444
$ws.Cells.Item(9, 1).Value = 20.66363120217269
$ws.Cells.Item(9, 2).Value = 17.51956762470492
$ws.Cells.Item(9, 3).Value = 24.24021904643489
$ws.Cells.Item(10, 1).Value = 1.345652398441431
$ws.Cells.Item(10, 2).Value = 1.062001570980459
$ws.Cells.Item(10, 3).Value = 1.834777650407574
$ws.Cells.Item(11, 1).Value = 7.779608989906619
$ws.Cells.Item(11, 2).Value = 5.24554499016497
$ws.Cells.Item(11, 3).Value = 12.00400399590753
$ws.Cells.Item(12, 1).Value = 21.72542333256477
$ws.Cells.Item(12, 2).Value = 18.81504466080679
$ws.Cells.Item(12, 3).Value = 25.69784968048811
$ws.Cells.Item(13, 1).Value = 2.090620064691609
$ws.Cells.Item(13, 2).Value = 1.493238504029998
$ws.Cells.Item(13, 3).Value = 3.160184548130703
$ws.Cells.Item(14, 1).Value = 3.902163322443779
$ws.Cells.Item(14, 2).Value = 2.603895573143855
$ws.Cells.Item(14, 3).Value = 6.207759772317279
$ws.Cells.Item(15, 1).Value = 25.60594210979114
$ws.Cells.Item(15, 2).Value = 11.01683556207869
$ws.Cells.Item(15, 3).Value = 40.5011045314051
$ws.Cells.Item(16, 1).Value = 22.50443637705029
$ws.Cells.Item(16, 2).Value = 18.87029079788068
$ws.Cells.Item(16, 3).Value = 26.20459653840652
$ws.Cells.Item(17, 1).Value = 21.78779876928554
$ws.Cells.Item(17, 2).Value = 18.54984337574082
$ws.Cells.Item(17, 3).Value = 25.55109723468735
$ws.Cells.Item(18, 1).Value = 17.62925599830132
$ws.Cells.Item(18, 2).Value = 12.57160958715028
$ws.Cells.Item(18, 3).Value = 23.26916816943996
$ws.Cells.Item(19, 1).Value = 7.872630629040973
$ws.Cells.Item(19, 2).Value = 5.28166901737814
$ws.Cells.Item(19, 3).Value = 12.13621763062513
$ws.Cells.Item(20, 1).Value = 25.4462726469938
$ws.Cells.Item(20, 2).Value = 21.98213128550441
$ws.Cells.Item(20, 3).Value = 29.49912199666365
$ws.Cells.Item(21, 1).Value = 25.23364652266176
$ws.Cells.Item(21, 2).Value = 18.01672490502324
$ws.Cells.Item(21, 3).Value = 32.06595100222099
$ws.Cells.Item(22, 1).Value = 1.345187281471877
$ws.Cells.Item(22, 2).Value = 1.058373854730492
$ws.Cells.Item(22, 3).Value = 1.87112430654644
$ws.Cells.Item(23, 1).Value = 28.1054204019949
$ws.Cells.Item(23, 2).Value = 12.98775318212981
$ws.Cells.Item(23, 3).Value = 43.85076428274915
$ws.Cells.Item(24, 1).Value = 22.50080727217914
$ws.Cells.Item(24, 2).Value = 18.873905632994
$ws.Cells.Item(24, 3).Value = 26.19407500130826
$ws.Cells.Item(25, 1).Value = 18.78897989754048
$ws.Cells.Item(25, 2).Value = 15.82913510761607
$ws.Cells.Item(25, 3).Value = 21.97026972366085
$ws.Cells.Item(26, 1).Value = 24.81522082281263
$ws.Cells.Item(26, 2).Value = 17.63452385514087
$ws.Cells.Item(26, 3).Value = 31.72316176153529
$ws.Cells.Item(27, 1).Value = 2.089879631621546
$ws.Cells.Item(27, 2).Value = 1.495554398578567
$ws.Cells.Item(27, 3).Value = 3.200912975293255
$ws.Cells.Item(28, 1).Value = 13.43450009631877
$ws.Cells.Item(28, 2).Value = 10.57551984512362
$ws.Cells.Item(28, 3).Value = 16.9853733310401
$ws.Cells.Item(29, 1).Value = 29.49280350105835
$ws.Cells.Item(29, 2).Value = 13.23565188293436
$ws.Cells.Item(29, 3).Value = 46.1371491085018
$ws.Cells.Item(30, 1).Value = 20.55226404532578
$ws.Cells.Item(30, 2).Value = 17.56959790377127
$ws.Cells.Item(30, 3).Value = 23.94800517646402
$ws.Cells.Item(31, 1).Value = 3.976440569913307
$ws.Cells.Item(31, 2).Value = 2.609044088049878
$ws.Cells.Item(31, 3).Value = 6.220159884583693
$ws.Cells.Item(32, 1).Value = 2.093725944236646
$ws.Cells.Item(32, 2).Value = 1.512787829427939
$ws.Cells.Item(32, 3).Value = 3.163517041945289
$ws.Cells.Item(33, 1).Value = 3.921209445194216
$ws.Cells.Item(33, 2).Value = 2.599515276370651
$ws.Cells.Item(33, 3).Value = 6.236237545554384
$ws.Cells.Item(34, 1).Value = 15.36330592822387
$ws.Cells.Item(34, 2).Value = 11.49711718232301
$ws.Cells.Item(34, 3).Value = 20.53613353911573
$ws.Cells.Item(35, 1).Value = 23.79795549470442
$ws.Cells.Item(35, 2).Value = 20.09485951279818
$ws.Cells.Item(35, 3).Value = 27.76568276076145
$ws.Cells.Item(36, 1).Value = 12.61298705037775
$ws.Cells.Item(36, 2).Value = 10.05065791668104
$ws.Cells.Item(36, 3).Value = 16.12635384807147
$ws.Cells.Item(37, 1).Value = 27.24053119280253
$ws.Cells.Item(37, 2).Value = 20.02984255167112
$ws.Cells.Item(37, 3).Value = 34.1755360770379
$ws.Cells.Item(38, 1).Value = 3.910077893354833
$ws.Cells.Item(38, 2).Value = 2.593569065320459
$ws.Cells.Item(38, 3).Value = 6.2666429407441
$ws.Cells.Item(39, 1).Value = 22.7292968763431
$ws.Cells.Item(39, 2).Value = 19.56781892300916
$ws.Cells.Item(39, 3).Value = 26.51215630817335
$ws.Cells.Item(40, 1).Value = 21.46648129787013
$ws.Cells.Item(40, 2).Value = 16.0439991624402
$ws.Cells.Item(40, 3).Value = 27.23771872120865
$ws.Cells.Item(41, 1).Value = 12.01468994799287
$ws.Cells.Item(41, 2).Value = 9.857912221068119
$ws.Cells.Item(41, 3).Value = 14.96182069475434
$ws.Cells.Item(42, 1).Value = 7.827396309315415
$ws.Cells.Item(42, 2).Value = 5.252742541167156
$ws.Cells.Item(42, 3).Value = 12.27555642366491
$ws.Cells.Item(43, 1).Value = 20.55226404532578
$ws.Cells.Item(43, 2).Value = 17.56959790377127
$ws.Cells.Item(43, 3).Value = 23.94800517646402
$ws.Cells.Item(44, 1).Value = 7.324651750261172
$ws.Cells.Item(44, 2).Value = 5.268188125513968
$ws.Cells.Item(44, 3).Value = 10.58042317634894
$ws.Cells.Item(45, 1).Value = 8.002011231768117
$ws.Cells.Item(45, 2).Value = 5.295585378083115
$ws.Cells.Item(45, 3).Value = 12.74283946839329
$ws.Cells.Item(46, 1).Value = 3.899261168214233
$ws.Cells.Item(46, 2).Value = 2.607815885709016
$ws.Cells.Item(46, 3).Value = 6.048136392168704
$ws.Cells.Item(47, 1).Value = 1.364013415023436
$ws.Cells.Item(47, 2).Value = 1.055614247508611
$ws.Cells.Item(47, 3).Value = 1.851319295429734
$ws.Cells.Item(48, 1).Value = 21.77381133130935
$ws.Cells.Item(48, 2).Value = 18.82652649674428
$ws.Cells.Item(48, 3).Value = 25.7611058011361
$ws.Cells.Item(49, 1).Value = 19.61257597255933
$ws.Cells.Item(49, 2).Value = 13.73246640420222
$ws.Cells.Item(49, 3).Value = 25.50161569634258
$ws.Cells.Item(50, 1).Value = 27.1917568931443
$ws.Cells.Item(50, 2).Value = 19.88245245681075
$ws.Cells.Item(50, 3).Value = 34.1529024587131
$ws.Cells.Item(51, 1).Value = 18.42879784190137
$ws.Cells.Item(51, 2).Value = 6.959129528158291
$ws.Cells.Item(51, 3).Value = 29.84979114662602
$ws.Cells.Item(52, 1).Value = 3.774785348618193
$ws.Cells.Item(52, 2).Value = 2.596981421510345
$ws.Cells.Item(52, 3).Value = 6.079559494950883
$ws.Cells.Item(53, 1).Value = 22.76280699234251
$ws.Cells.Item(53, 2).Value = 19.41980367754787
$ws.Cells.Item(53, 3).Value = 26.63180066305698
$ws.Cells.Item(54, 1).Value = 25.63570728911887
$ws.Cells.Item(54, 2).Value = 11.09084773863882
$ws.Cells.Item(54, 3).Value = 40.54069634484281
$ws.Cells.Item(55, 1).Value = 24.94566584941417
$ws.Cells.Item(55, 2).Value = 17.88908710530983
$ws.Cells.Item(55, 3).Value = 31.91573396106226
$ws.Cells.Item(56, 1).Value = 3.898376678198692
$ws.Cells.Item(56, 2).Value = 2.613055846654043
$ws.Cells.Item(56, 3).Value = 6.158106849405281
$ws.Cells.Item(57, 1).Value = 2.117278457128078
$ws.Cells.Item(57, 2).Value = 1.515095877654378
$ws.Cells.Item(57, 3).Value = 3.17020147291823
$ws.Cells.Item(58, 1).Value = 27.1843270231798
$ws.Cells.Item(58, 2).Value = 20.21363764448143
$ws.Cells.Item(58, 3).Value = 34.35849439711836
$ws.Cells.Item(59, 1).Value = 7.659139727458181
$ws.Cells.Item(59, 2).Value = 5.216199637170454
$ws.Cells.Item(59, 3).Value = 11.53600112225452
$ws.Cells.Item(60, 1).Value = 28.12849977053046
$ws.Cells.Item(60, 2).Value = 13.20192482061311
$ws.Cells.Item(60, 3).Value = 45.29619908247507
$ws.Cells.Item(61, 1).Value = 3.882043495823602
$ws.Cells.Item(61, 2).Value = 2.594177646257118
$ws.Cells.Item(61, 3).Value = 6.209182337654563
$ws.Cells.Item(62, 1).Value = 7.780200012641967
$ws.Cells.Item(62, 2).Value = 5.247999613055492
$ws.Cells.Item(62, 3).Value = 12.04164743485302
$ws.Cells.Item(63, 1).Value = 3.929000599989743
$ws.Cells.Item(63, 2).Value = 2.623880875568434
$ws.Cells.Item(63, 3).Value = 6.328636154866857
$ws.Cells.Item(64, 1).Value = 28.09404795275958
$ws.Cells.Item(64, 2).Value = 12.86086081213301
$ws.Cells.Item(64, 3).Value = 44.62712031535178
$ws.Cells.Item(65, 1).Value = 2.12772085338439
$ws.Cells.Item(65, 2).Value = 1.516397051789074
$ws.Cells.Item(65, 3).Value = 3.235053636391319
$ws.Cells.Item(66, 1).Value = 14.84213628682425
$ws.Cells.Item(66, 2).Value = 11.10741042114251
$ws.Cells.Item(66, 3).Value = 20.14497610595064
$ws.Cells.Item(67, 1).Value = 14.44759486070826
$ws.Cells.Item(67, 2).Value = 11.01382812680418
$ws.Cells.Item(67, 3).Value = 18.90648075643981
$ws.Cells.Item(68, 1).Value = 23.77649010148448
$ws.Cells.Item(68, 2).Value = 16.77154960584606
$ws.Cells.Item(68, 3).Value = 30.54449383645341
$ws.Cells.Item(69, 1).Value = 1.337221548398764
$ws.Cells.Item(69, 2).Value = 1.060647111690126
$ws.Cells.Item(69, 3).Value = 1.865484611650596
$ws.Cells.Item(70, 1).Value = 21.26947322133706
$ws.Cells.Item(70, 2).Value = 16.14589042436712
$ws.Cells.Item(70, 3).Value = 27.3648533074205
$ws.Cells.Item(71, 1).Value = 15.35006998807355
$ws.Cells.Item(71, 2).Value = 13.65688046764606
$ws.Cells.Item(71, 3).Value = 17.76119063550463
$ws.Cells.Item(72, 1).Value = 18.25335653799225
$ws.Cells.Item(72, 2).Value = 16.00621622659151
$ws.Cells.Item(72, 3).Value = 21.13098077222388

Write-Output "Updated 213 cells (A2:C72)"
